$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-4)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Update column B values (rows 3-4)
$ws.Range("B3").Value = 16
$ws.Range("B4").Value = 15

# Delete row 6 entirely (shift cells up), since row5 stays and row6 is removed
$ws.Range("A6:B6").Delete()
